# Generate Report for Handoff
# Updates the localization-status workbook with a new handoff file id
# (8ac53b9b-0631-4105-8595-710dacf8b319 replacing e5208c76-0d51-49c1-8e73-48d0b96689f1)
# and refreshed handoff/handback timestamps + content hash in the xlf file names.

$wb = $excel.ActiveWorkbook

$newId = "8ac53b9b-0631-4105-8595-710dacf8b319"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("G2").Value = "2016-09-03 05:03:19"

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newId.md"
}

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newId.md"
$wsZhCn.Range("G2").Value = "$newId.c3618a0c5960ef73649d25f8c0b434df5024cd11.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-03 05:03:14"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = "$newId.md"
}

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newId.md"
$wsDeDe.Range("G2").Value = "$newId.c3618a0c5960ef73649d25f8c0b434df5024cd11.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-03 05:03:19"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = "$newId.md"
}
